# Insert a new data row at row 119 (pushes existing rows 119:221 down to 120:222)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(119).Insert()

$r = 119
$ws.Cells.Item($r, 1).Value  = 3
$ws.Cells.Item($r, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item($r, 3).Value  = "Coquimbo"
$ws.Cells.Item($r, 4).Value  = 44827
$ws.Cells.Item($r, 5).Value  = 5
$ws.Cells.Item($r, 6).Value  = 100112010
$ws.Cells.Item($r, 7).Value  = "Achicoria"
$ws.Cells.Item($r, 8).Value  = "Sin especificar"
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 110
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 6500
$ws.Cells.Item($r, 13).Value = 6227
$ws.Cells.Item($r, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 16).Value = 389
$ws.Cells.Item($r, 17).Value = 16
$ws.Cells.Item($r, 18).Value = "Hortaliza"
